$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Village-name correction: "Y.Othakkadai" -> "Y.Othakadai"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Y.Othakkadai", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Y.Othakadai", 2)

# ---------------------------------------------------------------------------
# 2) "To limit the over absorption of public water source  from large
#    consumers." -> same sentence with "source" pluralised to "sources",
#    splitting the original single run into three runs.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("To limit the over absorption of public water source  from large consumers.")
if ($found) {
    $start = $r.Start

    $partA = "To limit the over absorption of public water "
    $partB = "source"
    $partC = "  from large consumers."

    # Replace "source" with "sources" in place.
    $rB = $d.Range($start + $partA.Length, $start + $partA.Length + $partB.Length)
    $rB.Text = "sources"
    $partB = "sources"

    # Re-derive the three sub-ranges over the (now longer) text and nudge
    # each one's formatting so the run splits land exactly on the word
    # boundaries, mirroring how the source document separates the runs.
    $off = $start
    $rA = $d.Range($off, $off + $partA.Length)
    $off = $off + $partA.Length
    $rB = $d.Range($off, $off + $partB.Length)
    $off = $off + $partB.Length
    $rC = $d.Range($off, $off + $partC.Length)

    $rA.Font.Superscript = $true
    $rA.Font.Superscript = $false
    $rB.Font.Superscript = $true
    $rB.Font.Superscript = $false
    $rC.Font.Superscript = $true
    $rC.Font.Superscript = $false
}

# ---------------------------------------------------------------------------
# 3) "Attended IIPC sponsored one day seminar on " ->
#    "Participated in IIPC sponsored a one day seminar on ", split across
#    five runs.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("Attended IIPC sponsored one day seminar on ")
if ($found2) {
    $start2 = $r2.Start
    $r2.Text = "Participated in IIPC sponsored a one day seminar on "

    $p1 = "Participated in "
    $p2 = "IIPC"
    $p3 = " "
    $p4 = "sponsored a one"
    $p5 = " day seminar on "

    $off2 = $start2
    $rA2 = $d.Range($off2, $off2 + $p1.Length)
    $off2 = $off2 + $p1.Length
    $rB2 = $d.Range($off2, $off2 + $p2.Length)
    $off2 = $off2 + $p2.Length
    $rC2 = $d.Range($off2, $off2 + $p3.Length)
    $off2 = $off2 + $p3.Length
    $rD2 = $d.Range($off2, $off2 + $p4.Length)
    $off2 = $off2 + $p4.Length
    $rE2 = $d.Range($off2, $off2 + $p5.Length)

    $rA2.Font.Superscript = $true
    $rA2.Font.Superscript = $false
    $rB2.Font.Superscript = $true
    $rB2.Font.Superscript = $false
    $rC2.Font.Superscript = $true
    $rC2.Font.Superscript = $false
    $rD2.Font.Superscript = $true
    $rD2.Font.Superscript = $false
    $rE2.Font.Superscript = $true
    $rE2.Font.Superscript = $false
}
